$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.830.21'
$ws.Range("E2").Value = '  -0.21%  '

$ws.Range("D3").Value = '3.144.84'
$ws.Range("E3").Value = '  -0.07%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").Value = '584.96'
$ws.Range("E5").Value = '  -0.97%  '

$ws.Range("D6").Value = '146.08'
$ws.Range("E6").Value = '  -1.02%  '

$ws.Range("E7").Value = '  +0.07%  '

$ws.Range("D8").Value = '3.131.13'
$ws.Range("E8").Value = '  -0.32%  '

$ws.Range("D9").Value = '0.528'
$ws.Range("E9").Value = '  -1.53%  '

$ws.Range("D10").Value = '0.161'
$ws.Range("E10").Value = '  +0.51%  '

$ws.Range("E11").Value = '  +0.59%  '

$ws.Range("D12").Value = '0.457'
$ws.Range("E12").Value = '  -2.92%  '

$ws.Range("D13").Value = '0.0000246'
$ws.Range("E13").Value = '  -3.00%  '

$ws.Range("D14").Value = '36.88'
$ws.Range("E14").Value = '  +2.44%  '

$ws.Range("D15").Value = '3.664.80'
$ws.Range("E15").Value = '  -0.09%  '

$ws.Range("E16").Value = '  -1.57%  '

$ws.Range("D17").Value = '63.636.10'
$ws.Range("E17").Value = '  -0.45%  '

$ws.Range("D18").Value = '3.141.38'
$ws.Range("E18").Value = '  -0.08%  '

$ws.Range("D19").Value = '7.06'
$ws.Range("E19").Value = '  -1.64%  '

$ws.Range("D20").Value = '463.69'
$ws.Range("E20").Value = '  -1.67%  '

$ws.Range("D21").Value = '14.27'
$ws.Range("E21").Value = '  +0.21%  '

$ws.Range("E22").Value = '  +0.11%  '

$ws.Range("D23").Value = '7.41'
$ws.Range("E23").Value = '  -2.07%  '

$ws.Range("D24").Value = '12.90'
$ws.Range("E24").Value = '  -3.68%  '

$ws.Range("D25").Value = '80.99'
$ws.Range("E25").Value = '  -1.94%  '

$ws.Range("E26").Value = '  +1.00%  '

$ws.Range("E27").Value = '  -0.13%  '

$ws.Range("D28").Value = '9.29'
$ws.Range("E28").Value = '  +6.17%  '

$ws.Range("D29").Value = '2.67'
$ws.Range("E29").Value = '  -1.17%  '

$ws.Range("B30").Value = 'FirstDigitalUSD'
$ws.Range("C30").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  -0.05%  '

$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").Value = '2.20'
$ws.Range("E31").Value = '  -0.61%  '

$ws.Range("D32").Value = '7.10'
$ws.Range("E32").Value = '  +4.05%  '

$ws.Range("D33").Value = '26.85'
$ws.Range("E33").Value = '  -0.92%  '

$ws.Range("E34").Value = '  -0.30%  '

$ws.Range("E35").Value = '  -3.61%  '

$ws.Range("E36").Value = '  -0.79%  '

$ws.Range("D37").Value = '2.31'
$ws.Range("E37").Value = '  -4.19%  '

$ws.Range("D38").Value = '3.31'
$ws.Range("E38").Value = '  -3.25%  '

$ws.Range("D39").Value = '6.00'
$ws.Range("E39").Value = '  -2.59%  '

$ws.Range("D40").Value = '51.22'

$ws.Range("D41").Value = '435.71'
$ws.Range("E41").Value = '  -3.05%  '

$ws.Range("D42").Value = '8.87'
$ws.Range("E42").Value = '  +1.47%  '

$ws.Range("D43").Value = '2.917.29'
$ws.Range("E43").Value = '  -0.35%  '

$ws.Range("D44").Value = '0.0370'
$ws.Range("E44").Value = '  -1.48%  '

$ws.Range("E45").Value = '  -1.64%  '

$ws.Range("D46").Value = '0.107'
$ws.Range("E46").Value = '  -4.51%  '

$ws.Range("D47").Value = '37.45'
$ws.Range("E47").Value = '  +7.29%  '

$ws.Range("D48").Value = '126.72'
$ws.Range("E48").Value = '  +2.17%  '

$ws.Range("E49").Value = '  +0.01%  '

$ws.Range("E50").Value = '  -1.60%  '

$ws.Range("D51").Value = '24.08'
$ws.Range("E51").Value = '  -3.32%  '
